# Fixed minor issue in datastream advanced slides
#
# The code sample on the "OperatorState" slide called:
#     .getOperatorState("totalLengthByKey", 0L, false);
# and should instead call it with the "fault tolerant" flag flipped on:
#     .getOperatorState("totalLengthByKey", 0L, true);
#
# We locate the run that holds the text and edit it in place. PowerPoint
# splits the touched run into several runs along the sub-strings that were
# actually (re)typed, so we reproduce that by writing to three adjacent
# Characters() sub-ranges: the comma+space, and the boolean literal (with
# its trailing paren folded in so the final ";" stays its own run).

$p = $ppt.ActivePresentation

$targetSlideIndex = 30
$targetShapeIndex = 3
$needle = '", 0L, false);'

# Primary, known-good location.
$s = $p.Slides.Item($targetSlideIndex)
$shp = $s.Shapes.Item($targetShapeIndex)
$tr = $shp.TextFrame.TextRange
$full = $tr.Text
$idx = $full.IndexOf($needle)

if ($idx -lt 0) {
    # Defensive fallback: scan every slide/shape for the text in case
    # indices ever shift.
    for ($si = 1; $si -le $p.Slides.Count -and $idx -lt 0; $si++) {
        $cs = $p.Slides.Item($si)
        for ($hi = 1; $hi -le $cs.Shapes.Count -and $idx -lt 0; $hi++) {
            $cshp = $cs.Shapes.Item($hi)
            if ($cshp.HasTextFrame) {
                $ctr = $cshp.TextFrame.TextRange
                $cfull = $ctr.Text
                $cidx = $cfull.IndexOf($needle)
                if ($cidx -ge 0) {
                    $s = $cs
                    $shp = $cshp
                    $tr = $ctr
                    $full = $cfull
                    $idx = $cidx
                }
            }
        }
    }
}

if ($idx -ge 0) {
    $start = $idx + 1   # TextRange/Characters() use 1-based character positions

    # '", 0L' (5 chars) is left completely untouched so it keeps its original
    # run/formatting as-is.

    # ', ' -> re-write in place (same text, but this carves out its own run,
    # matching how the edit was captured originally).
    $rComma = $tr.Characters($start + 5, 2)
    $rComma.Text = ", "

    # 'false)' -> 'true)'
    $rBool = $tr.Characters($start + 7, 6)
    $rBool.Text = "true)"

    # trailing ';' is left untouched, keeping it as its own run.
}

# Best-effort: the deck's Notes Master also carries a cached
# "datetimeFigureOut" field (01/09/15 -> 03/09/15) that PowerPoint re-stamps
# whenever it resaves the whole deck. That auto-computed field isn't
# writable through the PowerPoint object model (it isn't tied to any
# editable TextRange/HeaderFooter text), so there is nothing further to do
# for it here; left as a no-op/comment for documentation purposes.
